$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2
$ws.Range("C2").Value = 104.4425746990862
$ws.Range("D2").Value = 1282.730868820775
$ws.Range("E2").Value = 14906.39371105165
$ws.Range("F2").Value = 1282.730718813274
$ws.Range("G2").Value = 44947.63800973798
$ws.Range("H2").Value = 1.10702453768404
$ws.Range("J2").Value = 18.49878234783836
$ws.Range("K2").Value = 0.6927723515641337
$ws.Range("M2").Value = 1387.17314350486
$ws.Range("N2").Value = 29421.8001272939
$ws.Range("O2").Value = 1387.17329351236
$ws.Range("P2").Value = 80804.4521471122
$ws.Range("Q2").Value = 0.9781771705224789
$ws.Range("S2").Value = 41.97131683209776
$ws.Range("T2").Value = 0.7554541260111384
$ws.Range("V2").Value = 1262.674645098945
# Row 3
$ws.Range("C3").Value = 104.9481986106382
$ws.Range("D3").Value = 1281.782461697167
$ws.Range("E3").Value = 19888.94515174412
$ws.Range("F3").Value = 1281.782411694667
$ws.Range("G3").Value = 59001.91880459722
$ws.Range("H3").Value = 1.071728699732787
$ws.Range("J3").Value = 27.54146420985196
$ws.Range("K3").Value = 0.7429136612205511
$ws.Range("L3").Value = $False
$ws.Range("M3").Value = 1386.730660307805
$ws.Range("N3").Value = 38315.77779261894
$ws.Range("O3").Value = 1386.730610305305
$ws.Range("P3").Value = 110782.833857773
$ws.Range("Q3").Value = 1.031879739104645
$ws.Range("S3").Value = 29.36096108578936
$ws.Range("T3").Value = 0.7630173464862076
$ws.Range("V3").Value = 1261.760056663317
# Row 4
$ws.Range("C4").Value = 104.5378114030855
$ws.Range("D4").Value = 1282.549917620856
$ws.Range("E4").Value = 26125.69411200816
$ws.Range("F4").Value = 1282.547067478349
$ws.Range("G4").Value = 77936.44893853911
$ws.Range("H4").Value = 1.0898720731134
$ws.Range("J4").Value = 47.50709129768414
$ws.Range("K4").Value = 0.7108862664864701
$ws.Range("L4").Value = $True
$ws.Range("M4").Value = 1387.084928883934
$ws.Range("N4").Value = 50877.91180974832
$ws.Range("O4").Value = 1387.084878881434
$ws.Range("P4").Value = 138781.5532367061
$ws.Range("Q4").Value = 0.9919617373911227
$ws.Range("S4").Value = 38.83076522581072
$ws.Range("T4").Value = 0.7139670186414619
$ws.Range("U4").Value = $False
$ws.Range("V4").Value = 1262.594382681717
# Row 5
$ws.Range("C5").Value = 104.5928468061729
$ws.Range("D5").Value = 1282.422381017461
$ws.Range("E5").Value = 27783.20867138372
$ws.Range("F5").Value = 1282.422231009961
$ws.Range("G5").Value = 82189.03412892266
$ws.Range("H5").Value = 1.086462528188054
$ws.Range("J5").Value = 33.2358927304521
$ws.Range("K5").Value = 0.6975662738323346
$ws.Range("M5").Value = 1387.015127818634
$ws.Range("N5").Value = 53834.87909182344
$ws.Range("O5").Value = 1387.015077816133
$ws.Range("P5").Value = 145258.6941514138
$ws.Range("Q5").Value = 0.9989210625014198
$ws.Range("S5").Value = 116.6068168691055
$ws.Range("T5").Value = 0.6816395068060388
$ws.Range("U5").Value = $True
$ws.Range("V5").Value = 1262.422318144928
# Row 6
$ws.Range("C6").Value = 104.6487359290154
$ws.Range("D6").Value = 1282.315324617347
$ws.Range("E6").Value = 28729.21125653609
$ws.Range("F6").Value = 1282.315374619847
$ws.Range("G6").Value = 87787.2242002834
$ws.Range("H6").Value = 1.082845221236768
$ws.Range("J6").Value = 57.28451394872647
$ws.Range("K6").Value = 0.7652797283714029
$ws.Range("L6").Value = $True
$ws.Range("M6").Value = 1386.964160551363
$ws.Range("N6").Value = 55424.39769154091
$ws.Range("O6").Value = 1386.964110548863
$ws.Range("P6").Value = 153558.3155044282
$ws.Range("Q6").Value = 1.003900504424702
$ws.Range("S6").Value = 45.39056612562766
$ws.Range("T6").Value = 0.7242722724176467
$ws.Range("V6").Value = 1262.358228537804
# Row 7
$ws.Range("C7").Value = 104.7006076909645
$ws.Range("D7").Value = 1282.221380320299
$ws.Range("E7").Value = 28849.95166788357
$ws.Range("F7").Value = 1282.221330317799
$ws.Range("G7").Value = 84972.86807868385
$ws.Range("H7").Value = 1.075852099807332
$ws.Range("J7").Value = 35.13569201322791
$ws.Range("K7").Value = 0.712122562219541
$ws.Range("M7").Value = 1386.921988011264
$ws.Range("N7").Value = 55482.68370309147
$ws.Range("O7").Value = 1386.921938008764
$ws.Range("P7").Value = 154793.7462892019
$ws.Range("Q7").Value = 1.010704127227788
$ws.Range("S7").Value = 49.55675165777669
$ws.Range("T7").Value = 0.7263464159974007
$ws.Range("V7").Value = 1262.172490665632
# Row 8
$ws.Range("C8").Value = 104.8871729540144
$ws.Range("D8").Value = 1281.888682685162
$ws.Range("E8").Value = 30295.70183096123
$ws.Range("F8").Value = 1281.888532677662
$ws.Range("G8").Value = 89463.15198056828
$ws.Range("H8").Value = 1.071076413346379
$ws.Range("J8").Value = 39.42431652554543
$ws.Range("K8").Value = 0.7332669497009476
$ws.Range("M8").Value = 1386.775755634176
$ws.Range("N8").Value = 58207.23499135145
$ws.Range("O8").Value = 1386.775705631676
$ws.Range("P8").Value = 166953.0284489367
$ws.Range("Q8").Value = 1.026371972833134
$ws.Range("S8").Value = 45.81480602016917
$ws.Range("T8").Value = 0.7560447810855629
$ws.Range("U8").Value = $False
$ws.Range("V8").Value = 1261.86722328511
# Row 9
$ws.Range("C9").Value = 104.80945106739
$ws.Range("D9").Value = 1282.031124675205
$ws.Range("E9").Value = 29858.49283846282
$ws.Range("F9").Value = 1282.0291745777
$ws.Range("G9").Value = 89472.55024240245
$ws.Range("H9").Value = 1.071501073374309
$ws.Range("J9").Value = 57.15751791713139
$ws.Range("K9").Value = 0.753990971408187
$ws.Range("L9").Value = $True
$ws.Range("M9").Value = 1386.83867564759
$ws.Range("N9").Value = 57197.89006521182
$ws.Range("O9").Value = 1386.83862564509
$ws.Range("P9").Value = 162241.5092842063
$ws.Range("Q9").Value = 1.021358859110689
$ws.Range("S9").Value = 46.31357694180411
$ws.Range("T9").Value = 0.7405643240705246
$ws.Range("U9").Value = $False
$ws.Range("V9").Value = 1262.174421049026
